$wb = $excel.ActiveWorkbook
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "two_by_two_PriceinDem"

$ws.Range("B1").Value = "'benchmark"
$ws.Range("C1").Value = "'RA=157"
$ws.Range("D1").Value = "'eRA=.6"
$ws.Range("E1").Value = "'pr_Ud=3"
$ws.Range("F1").Value = "'prU2,eRA.5"
$ws.Range("G1").Value = "'prU.5,eRA.6"
$ws.Range("H1").Value = "'Itax=0.1"
$ws.Range("I1").Value = "'Otax=0.1"
$ws.Range("A2").Value = "'X"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1.0363877164248041
$ws.Range("D2").Value = 1.0365191201635342
$ws.Range("E2").Value = 1.0750070482317526
$ws.Range("F2").Value = 1.0595609099982375
$ws.Range("G2").Value = 0.99833470833488758
$ws.Range("H2").Value = 0.98248710709607356
$ws.Range("I2").Value = 0.94048657015713999
$ws.Range("A3").Value = "'Y"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1.0432700717660239
$ws.Range("D3").Value = 1.0430741066727305
$ws.Range("E3").Value = 0.9856633946311647
$ws.Range("F3").Value = 1.0087067867539614
$ws.Range("G3").Value = 1.1000067470017649
$ws.Range("H3").Value = 1.1229317376964334
$ws.Range("I3").Value = 1.1855074833626029
$ws.Range("A4").Value = "'U"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.0388246300487494
$ws.Range("D4").Value = 1.0389569047712837
$ws.Range("E4").Value = 1.077706227526201
$ws.Range("F4").Value = 1.062153685062508
$ws.Range("G4").Value = 1.0005254867585018
$ws.Range("H4").Value = 0.99843756359281455
$ws.Range("I4").Value = 0.99193526198210025
$ws.Range("A5").Value = "'PX"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1.0023513532505028
$ws.Range("D5").Value = 1.0023518955237618
$ws.Range("E5").Value = 1.0025108479973333
$ws.Range("F5").Value = 1.0024470278582356
$ws.Range("G5").Value = 1.0021944327924277
$ws.Range("H5").Value = 1.0162347743614697
$ws.Range("I5").Value = 1.0547043343918923
$ws.Range("A6").Value = "'PY"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.99573893485424703
$ws.Range("D6").Value = 0.9957379554075485
$ws.Range("E6").Value = 0.99545092202416496
$ws.Range("F6").Value = 0.99556615181893504
$ws.Range("G6").Value = 0.99602242515110884
$ws.Range("H6").Value = 0.97114380037066661
$ws.Range("I6").Value = 0.90770374205596738
$ws.Range("A7").Value = "'PU"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("A8").Value = "'PL"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.94438602708152408
$ws.Range("D8").Value = 0.94437357942991329
$ws.Range("E8").Value = 0.94073225025073981
$ws.Range("F8").Value = 0.94219250114330766
$ws.Range("G8").Value = 0.94799525094925208
$ws.Range("H8").Value = 0.89903690333116504
$ws.Range("I8").Value = 0.84467427236056347
$ws.Range("A9").Value = "'PK"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1.0388246300416317
$ws.Range("D9").Value = 1.0388337447788984
$ws.Range("E9").Value = 1.0415087178322844
$ws.Range("F9").Value = 1.0404339284150994
$ws.Range("G9").Value = 1.0361901997038321
$ws.Range("H9").Value = 1.0329715504842174
$ws.Range("I9").Value = 0.96149717694622538
$ws.Range("A10").Value = "'SX"
$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 80
$ws.Range("D10").Value = 80
$ws.Range("E10").Value = 80
$ws.Range("F10").Value = 80
$ws.Range("G10").Value = 80
$ws.Range("H10").Value = 80
$ws.Range("I10").Value = 80
$ws.Range("A11").Value = "'SY"
$ws.Range("B11").Value = 54
$ws.Range("C11").Value = 54
$ws.Range("D11").Value = 54
$ws.Range("E11").Value = 54
$ws.Range("F11").Value = 54
$ws.Range("G11").Value = 54
$ws.Range("H11").Value = 54
$ws.Range("I11").Value = 53.999999999999993
$ws.Range("A12").Value = "'SU"
$ws.Range("B12").Value = 124
$ws.Range("C12").Value = 124
$ws.Range("D12").Value = 124
$ws.Range("E12").Value = 124
$ws.Range("F12").Value = 124.00000000000001
$ws.Range("G12").Value = 124
$ws.Range("H12").Value = 124
$ws.Range("I12").Value = 124
$ws.Range("A13").Value = "'DXL"
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 31.841365430850797
$ws.Range("D13").Value = 31.841802353436101
$ws.Range("E13").Value = 31.970122669751873
$ws.Range("F13").Value = 31.918541910760322
$ws.Range("G13").Value = 31.715172574535519
$ws.Range("H13").Value = 30.827982414866039
$ws.Range("I13").Value = 30.648734246825121
$ws.Range("A14").Value = "'DXK"
$ws.Range("B14").Value = 50
$ws.Range("C14").Value = 48.244493065345502
$ws.Range("D14").Value = 48.244095867355547
$ws.Range("E14").Value = 48.127818367352575
$ws.Range("F14").Value = 48.174468386726872
$ws.Range("G14").Value = 48.359578824249922
$ws.Range("H14").Value = 49.189872358276311
$ws.Range("I14").Value = 49.362282267302113
$ws.Range("A15").Value = "'DYL"
$ws.Range("B15").Value = 24
$ws.Range("C15").Value = 25.305048730086334
$ws.Range("D15").Value = 25.305357380693888
$ws.Range("E15").Value = 25.395985012953503
$ws.Range("F15").Value = 25.359560402635765
$ws.Range("G15").Value = 25.215883918854793
$ws.Range("H15").Value = 25.924910448641455
$ws.Range("I15").Value = 25.790876462309562
$ws.Range("A16").Value = "'DYK"
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = 28.755737186305545
$ws.Range("D16").Value = 28.755456597924343
$ws.Range("E16").Value = 28.673334317227702
$ws.Range("F16").Value = 28.706277005080288
$ws.Range("G16").Value = 28.837054010926813
$ws.Range("H16").Value = 28.204372131496193
$ws.Range("I16").Value = 28.321572766511746
$ws.Range("A17").Value = "'DUX"
$ws.Range("B17").Value = 80
$ws.Range("C17").Value = 79.812333011330708
$ws.Range("D17").Value = 79.812289832779825
$ws.Range("E17").Value = 79.799635245655324
$ws.Range("F17").Value = 79.804715637616169
$ws.Range("G17").Value = 79.824829775888887
$ws.Range("H17").Value = 78.721966634399237
$ws.Range("I17").Value = 75.850641162032716
$ws.Range("A18").Value = "'DUY"
$ws.Range("B18").Value = 44
$ws.Range("C18").Value = 44.188289178863421
$ws.Range("D18").Value = 44.188332644183888
$ws.Range("E18").Value = 44.201074132845875
$ws.Range("F18").Value = 44.195958168736873
$ws.Range("G18").Value = 44.175712201785153
$ws.Range("H18").Value = 45.307399360636381
$ws.Range("I18").Value = 48.473965635901315
$ws.Range("A19").Value = "'RA"
$ws.Range("B19").Value = 134
$ws.Range("C19").Value = 139.20250041197306
$ws.Range("D19").Value = 139.20249020044872
$ws.Range("E19").Value = 139.20019309147671
$ws.Range("F19").Value = 139.20094884112044
$ws.Range("G19").Value = 139.20613388269214
$ws.Range("H19").Value = 138.76352763367868
$ws.Range("I19").Value = 137.46365821827374
$ws.Range("A20").Value = "'DU"
$ws.Range("B20").Value = 124
$ws.Range("C20").Value = 128.81425411257209
$ws.Range("D20").Value = 128.83065617550108
$ws.Range("E20").Value = 133.63557221313272
$ws.Range("F20").Value = 131.70705694774796
$ws.Range("G20").Value = 124.06516035730138
$ws.Range("H20").Value = 123.80625788538022
$ws.Range("I20").Value = 122.99997248577242
$ws.Range("A21").Value = "'DY"
$ws.Range("B21").Value = 10
$ws.Range("C21").Value = 10.432700716801408
$ws.Range("D21").Value = 10.416228455108469
$ws.Range("E21").Value = 5.5900504537469322
$ws.Range("F21").Value = 7.5272666509210788
$ws.Range("G21").Value = 15.201438384375418
$ws.Range("H21").Value = 15.401704405248296
$ws.Range("I21").Value = 15.934368299221513
$ws.Range("A22").Value = "'CWI"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 1.0391557242667655
$ws.Range("D22").Value = 1.0391554626179134
$ws.Range("E22").Value = 1.0548832266630734
$ws.Range("F22").Value = 1.0454956618532232
$ws.Range("G22").Value = 1.0534806670463777
$ws.Range("H22").Value = 1.0530112713727526
$ws.Range("I22").Value = 1.0506793942727906
$ws.Range("A23").Value = "'PX/PX"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("A24").Value = "'PY/PX"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 0.99340309326184617
$ws.Range("D24").Value = 0.99340157868134993
$ws.Range("E24").Value = 0.99295775603099801
$ws.Range("F24").Value = 0.9931359205543242
$ws.Range("G24").Value = 0.99384150675820293
$ws.Range("H24").Value = 0.9556293731247929
$ws.Range("I24").Value = 0.86062388525151878
$ws.Range("A25").Value = "'PU/PX"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 0.99765416264179452
$ws.Range("D25").Value = 0.99765362291001314
$ws.Range("E25").Value = 0.99749544057069406
$ws.Range("F25").Value = 0.9975589454702023
$ws.Range("G25").Value = 0.9978103721986229
$ws.Range("H25").Value = 0.98402458292999229
$ws.Range("I25").Value = 0.94813301452540921
$ws.Range("A26").Value = "'PL/PX"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 0.94217065105862896
$ws.Range("D26").Value = 0.94215772289875011
$ws.Range("E26").Value = 0.93837613042292212
$ws.Range("F26").Value = 0.93989255787045034
$ws.Range("G26").Value = 0.94591949419220012
$ws.Range("H26").Value = 0.88467441383912149
$ws.Range("I26").Value = 0.8008635641452776
$ws.Range("A27").Value = "'PK/PX"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 1.0363877164158561
$ws.Range("D27").Value = 1.036396249079844
$ws.Range("E27").Value = 1.0389001973523331
$ws.Range("F27").Value = 1.0378941724611865
$ws.Range("G27").Value = 1.0339213288350459
$ws.Range("H27").Value = 1.0164693991437794
$ws.Range("I27").Value = 0.91162721683569548
$ws.Range("A28").Value = "'RA/PX"
$ws.Range("B28").Value = 134
$ws.Range("C28").Value = 138.87595398615105
$ws.Range("D28").Value = 138.87586866657327
$ws.Range("E28").Value = 138.85155793530822
$ws.Range("F28").Value = 138.86115173439967
$ws.Range("G28").Value = 138.90132426182038
$ws.Range("H28").Value = 136.5467224056251
$ws.Range("I28").Value = 130.3338326541824
